$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 233, shifting existing rows 233:284 down to 234:285
$ws.Rows("233:233").Insert()

# Populate the newly inserted row 233 with the new data record
$ws.Range("A233").Value = 9
$ws.Range("B233").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C233").Value = "Metropolitana"
$ws.Range("D233").Value = 44508
$ws.Range("E233").Value = 13
$ws.Range("F233").Value = 100112031
$ws.Range("G233").Value = "Poroto verde"
$ws.Range("H233").Value = "Magnum"
$ws.Range("I233").Value = "Primera"
$ws.Range("J233").Value = 25
$ws.Range("K233").Value = 32000
$ws.Range("L233").Value = 34000
$ws.Range("M233").Value = 32960
$ws.Range("N233").Value = "$/malla 25 kilos"
$ws.Range("O233").Value = "Perú"
$ws.Range("P233").Value = 1318
$ws.Range("Q233").Value = 25
$ws.Range("R233").Value = "Hortaliza"
